$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Data updates: "Fatalities by Age Group" sheet
# ---------------------------------------------------------------------------
$wsFA = $wb.Worksheets.Item("Fatalities by Age Group")
$wsFA.Range("B5").Value = 296
$wsFA.Range("B6").Value = 963
$wsFA.Range("B7").Value = 2777
$wsFA.Range("B8").Value = 6131
$wsFA.Range("B9").Value = 5036
$wsFA.Range("B10").Value = 6430
$wsFA.Range("B11").Value = 7046
$wsFA.Range("B12").Value = 6931
$wsFA.Range("B13").Value = 17262
$wsFA.Range("B15").Formula = "=SUM(B2:B14)"

# ---------------------------------------------------------------------------
# Data updates: "Fatalities by Gender" sheet
# ---------------------------------------------------------------------------
$wsFG = $wb.Worksheets.Item("Fatalities by Gender")
$wsFG.Range("B2").Value = 22195
$wsFG.Range("B3").Value = 30735
$wsFG.Range("B5").Formula = "=SUM(B2:B4)"

# ---------------------------------------------------------------------------
# Data updates: "Fatalities by Race-Ethnicity" sheet
# ---------------------------------------------------------------------------
$wsFR = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$wsFR.Range("B2").Value = 1105
$wsFR.Range("B3").Value = 5387
$wsFR.Range("B4").Value = 24389
$wsFR.Range("B5").Value = 295
$wsFR.Range("B6").Value = 21732
$wsFR.Range("B8").Formula = "=SUM(B2:B7)"

# ---------------------------------------------------------------------------
# View-state updates (selection / active cell) per sheet, matching the
# author's saved workbook state.
# ---------------------------------------------------------------------------
$wsCA = $wb.Worksheets.Item("Cases by Age Group")
$wsCA.Range("E14").Select()

$wsCG = $wb.Worksheets.Item("Cases by Gender")
$wsCG.Range("B2:B4").Select()

$wsFA.Range("E9").Select()
$wsFG.Range("D12").Select()
$wsFR.Range("D9").Select()

$wsCR = $wb.Worksheets.Item("Cases by RaceEthnicity")
$wsCR.Range("I10").Select()

# The active sheet moved from "Cases by Age Group" to "Cases by RaceEthnicity"
# in the saved file (tabSelected moved accordingly).
$wsCR.Activate()
